# Corrección en parseo de Renta Fija y Venta Simultanea.
# - Para renta fija no hay que considerar los registros "Retrov Nominal"
#   (el antiguo BTP0600433 en la fila 15) -> se elimina esa fila.
# - Para la venta de simultaneas no se estaba haciendo el cálculo para
#   determinar el valor de la cantidad -> se recalcula la columna E
#   (cantidad) de la fila LTM/VENTA/SIMULTANEA.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "BTP0600433 / COMPRA / RENTA FIJA" row (row 15). Deleting the
# entire row shifts rows 16-20 up to 15-19 and shrinks the used range to
# A1:J19, matching the rest of the diff automatically.
$ws.Rows.Item(15).Delete()

# The simultanea "VENTA" row (now row 15, nemotecnico LTM) needs its
# "cantidad" (column E) recalculated - it was not being computed before.
$ws.Range("E15").Value = 12639064.464
